$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Updated values per row (B, C, D, E, F, H) - column G (Площадь) is unchanged
$data = @{
    2 = @(1264.4404296875, 0.9486, 0.9150000214576721, 1.483399987220764, 0.7125999927520752, 0.7645999999999999)
    3 = @(1286.101318359375, 1.0289, 0.9416, 2.034899950027466, 0.6187000274658203, 1.0003)
    4 = @(872.668212890625, 1.0195, 0.9415, 2.197400093078613, 0.7113000154495239, 0.9993)
    5 = @(831.2965087890625, 0.8769, 0.8716, 1.308200001716614, 0.5623999834060669, 0.3804)
    6 = @(1116.104736328125, 0.883, 0.8754, 1.193199992179871, 0.6590999960899353, 0.4135)
    7 = @(873.2382202148438, 0.8794, 0.8737000226974487, 1.094799995422363, 0.7120000123977661, 0.3985)
    8 = @(929.8483276367188, 0.8332000000000001, 0.8294, 1.138299942016602, 0.7023000121116638, 0.0064)
    9 = @(7173.69775390625, 0.9246, 0.8848, 2.197400093078613, 0.5623999834060669, 3.963)
}

foreach ($row in $data.Keys) {
    $vals = $data[$row]
    $ws.Cells.Item($row, 2).Value = $vals[0]
    $ws.Cells.Item($row, 3).Value = $vals[1]
    $ws.Cells.Item($row, 4).Value = $vals[2]
    $ws.Cells.Item($row, 5).Value = $vals[3]
    $ws.Cells.Item($row, 6).Value = $vals[4]
    $ws.Cells.Item($row, 8).Value = $vals[5]
}
